$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3857.5
$ws.Range("I15").Value = 3857.5
$ws.Range("K15").Value = 11572.5
$ws.Range("M15").Value = -11403.5
$ws.Range("H86").Value = 1874.5
$ws.Range("I86").Value = 897.2
$ws.Range("J86").Value = 2572.5715
$ws.Range("K86").Value = 897.2
$ws.Range("L86").Value = 2572.5715
$ws.Range("M86").Value = 225.8
$ws.Range("N86").Value = -4818.5715
$ws.Range("H89").Value = 1874.5
$ws.Range("I89").Value = 897.2
$ws.Range("J89").Value = 2572.5715
$ws.Range("K89").Value = 4486
$ws.Range("L89").Value = 12862.8575
$ws.Range("M89").Value = 1130
$ws.Range("N89").Value = -24094.8575
$ws.Range("H99").Value = 666
$ws.Range("I99").Value = 412
$ws.Range("J99").Value = 1237.5
$ws.Range("K99").Value = 1236
$ws.Range("L99").Value = 3712.5
$ws.Range("M99").Value = 262
$ws.Range("N99").Value = -6708.5
$ws.Range("H100").Value = 2278.5
$ws.Range("I100").Value = 2043.7858
$ws.Range("J100").Value = 3100
$ws.Range("K100").Value = 2043.7858
$ws.Range("L100").Value = 3100
$ws.Range("M100").Value = -1502.7858
$ws.Range("N100").Value = -4182
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1739.4286
$ws.Range("I2").Value = 1712.9474
$ws.Range("J2").Value = 1770.875
$ws.Range("K2").Value = 1712.9474
$ws.Range("L2").Value = 1770.875
$ws.Range("M2").Value = -1599.9474
$ws.Range("N2").Value = -1996.875
$ws.Range("H110").Value = 967.8333
$ws.Range("I110").Value = 802.5909
$ws.Range("J110").Value = 1422.25
$ws.Range("K110").Value = 802.5909
$ws.Range("L110").Value = 1422.25
$ws.Range("M110").Value = 1242.4091
$ws.Range("N110").Value = -5512.25
$ws.Range("H116").Value = 1739.4286
$ws.Range("I116").Value = 1712.9474
$ws.Range("J116").Value = 1770.875
$ws.Range("K116").Value = 1712.9474
$ws.Range("L116").Value = 1770.875
$ws.Range("M116").Value = 581.0526
$ws.Range("N116").Value = -6358.875
$ws.Range("H132").Value = 40685.094
$ws.Range("I132").Value = 28724.918
$ws.Range("J132").Value = 66716.06
$ws.Range("K132").Value = 86174.754
$ws.Range("L132").Value = 200148.18
$ws.Range("M132").Value = -83644.754
$ws.Range("N132").Value = -205208.18
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1739.4286
$ws.Range("I3").Value = 1712.9474
$ws.Range("J3").Value = 1770.875
$ws.Range("K3").Value = 1712.9474
$ws.Range("L3").Value = 1770.875
$ws.Range("M3").Value = -1598.9474
$ws.Range("N3").Value = -1998.875
$ws.Range("H99").Value = 1855.4615
$ws.Range("I99").Value = 1957.8889
$ws.Range("J99").Value = 1625
$ws.Range("K99").Value = 1957.8889
$ws.Range("L99").Value = 1625
$ws.Range("M99").Value = -459.8888999999999
$ws.Range("N99").Value = -4621
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5020.5
$ws.Range("I62").Value = 3876.25
$ws.Range("J62").Value = 5783.3335
$ws.Range("K62").Value = 3876.25
$ws.Range("L62").Value = 5783.3335
$ws.Range("M62").Value = -3252.25
$ws.Range("N62").Value = -7031.3335
$ws.Range("H65").Value = 5020.5
$ws.Range("I65").Value = 3876.25
$ws.Range("J65").Value = 5783.3335
$ws.Range("K65").Value = 19381.25
$ws.Range("L65").Value = 28916.6675
$ws.Range("M65").Value = -16261.25
$ws.Range("N65").Value = -35156.6675
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 939.53845
$ws.Range("J131").Value = 1015.94116
$ws.Range("L131").Value = 3047.82348
$ws.Range("N131").Value = -13127.82348
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58610.95
$ws.Range("I70").Value = 96246.17999999999
$ws.Range("J70").Value = 6862.5
$ws.Range("K70").Value = 96246.17999999999
$ws.Range("L70").Value = 6862.5
$ws.Range("M70").Value = -95976.17999999999
$ws.Range("N70").Value = -7402.5
$ws.Range("H73").Value = 58610.95
$ws.Range("I73").Value = 96246.17999999999
$ws.Range("J73").Value = 6862.5
$ws.Range("K73").Value = 96246.17999999999
$ws.Range("L73").Value = 6862.5
$ws.Range("M73").Value = -95310.17999999999
$ws.Range("N73").Value = -8734.5
$ws.Range("H80").Value = 4560
$ws.Range("I80").Value = 4600
$ws.Range("J80").Value = 4553.3335
$ws.Range("K80").Value = 4600
$ws.Range("L80").Value = 4553.3335
$ws.Range("M80").Value = -3602
$ws.Range("N80").Value = -6549.3335
$ws.Range("H83").Value = 4560
$ws.Range("I83").Value = 4600
$ws.Range("J83").Value = 4553.3335
$ws.Range("K83").Value = 23000
$ws.Range("L83").Value = 22766.6675
$ws.Range("M83").Value = -18008
$ws.Range("N83").Value = -32750.6675
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1781.8334
$ws.Range("I82").Value = 1345.5
$ws.Range("K82").Value = 1345.5
$ws.Range("M82").Value = -984.5
$ws.Range("H85").Value = 1781.8334
$ws.Range("I85").Value = 1345.5
$ws.Range("K85").Value = 1345.5
$ws.Range("M85").Value = -97.5
$ws.Range("H93").Value = 1563.1818
$ws.Range("I93").Value = 1466.6666
$ws.Range("J93").Value = 1997.5
$ws.Range("K93").Value = 1466.6666
$ws.Range("L93").Value = 1997.5
$ws.Range("M93").Value = -218.6666
$ws.Range("N93").Value = -4493.5
$ws.Range("H100").Value = 1476.2273
$ws.Range("J100").Value = 1973.75
$ws.Range("L100").Value = 1973.75
$ws.Range("N100").Value = -3055.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1903.262
$ws.Range("I81").Value = 781.3125
$ws.Range("J81").Value = 2593.6924
$ws.Range("K81").Value = 1562.625
$ws.Range("L81").Value = 5187.3848
$ws.Range("M81").Value = -501.625
$ws.Range("N81").Value = -7309.3848
$ws.Range("H84").Value = 1903.262
$ws.Range("I84").Value = 781.3125
$ws.Range("J84").Value = 2593.6924
$ws.Range("K84").Value = 7813.125
$ws.Range("L84").Value = 25936.924
$ws.Range("M84").Value = -2509.125
$ws.Range("N84").Value = -36544.924
$ws.Range("H96").Value = 1412.6923
$ws.Range("I96").Value = 1497.3334
$ws.Range("J96").Value = 1297.2727
$ws.Range("K96").Value = 1497.3334
$ws.Range("L96").Value = 1297.2727
$ws.Range("M96").Value = -124.3334
$ws.Range("N96").Value = -4043.2727
